$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.3048080303191223
$ws.Range("C2").Value = 10.29869402782916
$ws.Range("D2").Value = 0.8054896365839992
$ws.Range("E2").Value = 645.3272768299601
$ws.Range("G2").Value = 656.7362685246924

# Row 3
$ws.Range("B3").Value = [double]"3.996802888650564e-14"
$ws.Range("C3").Value = 0.3127903958511391
$ws.Range("D3").Value = 3.900430680208489
$ws.Range("E3").Value = 2367095152636972
$ws.Range("G3").Value = 2367095152636976

# Row 4
$ws.Range("B4").Value = 0.01514828764759746
$ws.Range("C4").Value = 1.667794583268128
$ws.Range("D4").Value = 26.21740644021617
$ws.Range("E4").Value = 645.3272768299601
$ws.Range("G4").Value = 673.2276261410919
